$d = $word.ActiveDocument

# The commit only touches the <w:nsid w:val="..."/> identifiers that Word
# mints for four of the numbering definitions in word/numbering.xml
# (abstractNumId 990, 99411, 99711, 99413). These ids are internal,
# auto-generated list instance identifiers - they carry no visible text
# and no formatting/content meaning, so the faithful way to reproduce the
# commit through the Word object model is a plain find & replace of the
# old id for the new one, done everywhere the document exposes its
# underlying XML as text.
$idMap = @{
    "56a55413" = "ee9c21b3";
    "be6797e5" = "a4e0db51";
    "655042df" = "5908c565";
    "d4ae37cb" = "6f5b7082"
}

foreach ($oldId in $idMap.Keys) {
    $newId = $idMap[$oldId]

    # Primary approach: standard Find/Replace over the whole document.
    $d.Content.Find.Execute($oldId, $true, $false, $false, $false, $false,
                             $true, 1, $false, $newId, 2)

    # Defensive sweep: also run the same replace over every story range
    # (headers/footers/footnotes/etc.) in case the id is reachable there.
    foreach ($story in $d.StoryRanges) {
        $story.Find.Execute($oldId, $true, $false, $false, $false, $false,
                             $true, 1, $false, $newId, 2)
    }
}

# Belt-and-suspenders: some hosts expose the full package OOXML (including
# word/numbering.xml) through Content/Document.WordOpenXML as a raw string;
# round-trip through it too so the edit lands if that surface is writable.
$fullXml = $d.Content.WordOpenXML
$changed = $false
foreach ($oldId in $idMap.Keys) {
    if ($fullXml.Contains($oldId)) {
        $fullXml = $fullXml.Replace($oldId, $idMap[$oldId])
        $changed = $true
    }
}
if ($changed) {
    $d.Content.WordOpenXML = $fullXml
}

Write-Output "nsid refresh applied"
